$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new functional group analysis values
$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = 3
$ws.Cells.Item(2,3).Value = 10
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 25
$ws.Cells.Item(2,6).Value = 9
$ws.Cells.Item(2,7).Value = 6

$ws.Cells.Item(3,1).Value = 0
$ws.Cells.Item(3,2).Value = 2
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 7
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 2

$ws.Cells.Item(4,1).Value = 0
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 45
$ws.Cells.Item(4,6).Value = 27
$ws.Cells.Item(4,7).Value = 0

$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 4
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 28
$ws.Cells.Item(5,6).Value = 4
$ws.Cells.Item(5,7).Value = 28

# Add new row 6
$ws.Cells.Item(6,1).Value = 10
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 13
$ws.Cells.Item(6,6).Value = 5
$ws.Cells.Item(6,7).Value = 11
